$d = $word.ActiveDocument

# 1. Fix typo "neearly" -> "nearly" in the Profile paragraph.
$d.Content.Find.Execute("neearly finished", $true, $false, $false, $false, $false,
                         $true, 1, $false, "nearly finished", 2) | Out-Null

# 2 & 3. Remove the two stray " <a" run pairs that precede the Github and
# Linkedin hyperlinks (leaving the following " " run and the hyperlink
# itself untouched). Using Find + Range.Delete (rather than a text
# Replace) keeps the hyperlink fields intact.
for ($i = 0; $i -lt 2; $i++) {
    $r = $d.Content
    $r.Find.Execute(" <a", $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
    if ($r.Find.Found) {
        $r.Delete()
    }
}

# 4. Add "further" before "utilized" in the Tampere paragraph.
$d.Content.Find.Execute("then utilized in a", $true, $false, $false, $false, $false,
                         $true, 1, $false, "then further utilized in a", 2) | Out-Null
